$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.504.51"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +5.63%  '
$ws.Range('D3').Value = "'1.819.97"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +6.01%  '
$ws.Range('D4').Value = "'1.004"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.25%  '
$ws.Range('D5').Value = "'345.66"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.61%  '
$ws.Range('D6').Value = "'1.000"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('D7').Value = "'0.3825"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.70%  '
$ws.Range('D8').Value = "'0.3525"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +5.87%  '
$ws.Range('D9').Value = "'49.50"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.98%  '
$ws.Range('D10').Value = "'1.239"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.78%  '
$ws.Range('D11').Value = "'0.07813"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.50%  '
$ws.Range('E12').Value = '  +0.31%  '
$ws.Range('D13').Value = "'22.28"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +10.99%  '
$ws.Range('D14').Value = "'6.650"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.37%  '
$ws.Range('D15').Value = "'7.266"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.15%  '
$ws.Range('D16').Value = "'1.817.79"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.12%  '
$ws.Range('D17').Value = "'0.00001129"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.86%  '
$ws.Range('D18').Value = "'0.06730"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.42%  '
$ws.Range('D19').Value = "'86.44"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.33%  '
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D21').Value = "'17.72"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +8.52%  '
$ws.Range('D22').Value = "'6.561"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.13%  '
$ws.Range('D23').Value = "'13.33"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.38%  '
$ws.Range('D24').Value = "'27.529.23"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.93%  '
$ws.Range('D25').Value = "'2.457"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('D26').Value = "'2.690"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.36%  '
$ws.Range('D27').Value = "'22.24"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +15.31%  '
$ws.Range('D28').Value = "'1.515"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +15.61%  '
$ws.Range('D29').Value = "'153.56"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.46%  '
$ws.Range('D30').Value = "'2.021.20"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.08%  '
$ws.Range('D31').Value = "'136.99"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.09%  '
$ws.Range('D32').Value = "'6.399"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.33%  '
$ws.Range('D33').Value = "'4.084"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.59%  '
$ws.Range('D34').Value = "'14.11"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +9.22%  '
$ws.Range('D35').Value = "'0.08791"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.10%  '
$ws.Range('D36').Value = "'1.700"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.17%  '
$ws.Range('D37').Value = "'5.680"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.10%  '
$ws.Range('D38').Value = "'0.7093"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +15.17%  '
$ws.Range('D39').Value = "'0.2289"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.49%  '
$ws.Range('D40').Value = "'0.06551"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.62%  '
$ws.Range('D41').Value = "'0.02426"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.19%  '
$ws.Range('D42').Value = "'9.031"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.08%  '
$ws.Range('D43').Value = "'1.301"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.64%  '
$ws.Range('D44').Value = "'14.83"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.00%  '
$ws.Range('D45').Value = "'0.6633"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +13.11%  '
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('D47').Value = "'3.979"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.93%  '
$ws.Range('D48').Value = "'2.196"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +9.44%  '
$ws.Range('D49').Value = "'133.03"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.52%  '
$ws.Range('D50').Value = "'0.07374"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.48%  '
$ws.Range('D51').Value = "'80.97"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.06%  '
